# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
# timestamps and sets the Priority column to "ht" for the rows that were
# just handed off.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$rows = @(7, 8, 9, 10, 12, 14)

foreach ($r in $rows) {
    # Overview sheet: column G = "Latest HO Xliff Generate Date"
    $overview.Range("G$r").Value = "2016-08-30 02:23:09"

    # zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime"
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-08-30 02:22:58"

    # de-de sheet: column E = "Priority", column H = "Latest Handoff Datetime"
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-08-30 02:23:09"
}
